$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 27 de Agosto de 2020 a las 15:54"

# Row 4
$ws.Range("B4").Value = 6003322
$ws.Range("C4").Value = 2957
$ws.Range("E4").Value = 2504959
$ws.Range("G4").Value = 46
$ws.Range("H4").Value = 183699

# Row 6
$ws.Range("B6").Value = 3333732
$ws.Range("C6").Value = 25983
$ws.Range("D6").Value = 2541192
$ws.Range("E6").Value = 731692
$ws.Range("G6").Value = 219
$ws.Range("H6").Value = 60848

# Row 14
$ws.Range("D14").Value = 274458
$ws.Range("E14").Value = 87786
$ws.Range("G14").Value = 105
$ws.Range("H14").Value = 7944

# Row 17
$ws.Range("B17").Value = 311855
$ws.Range("C17").Value = 1019
$ws.Range("D17").Value = 286255
$ws.Range("E17").Value = 21815
$ws.Range("G17").Value = 30
$ws.Range("H17").Value = 3785

# Row 23
$ws.Range("B23").Value = 239235
$ws.Range("C23").Value = 235
$ws.Range("E23").Value = 15648
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 9354

# Row 38
$ws.Range("A38").Value = "Oman"
$ws.Range("B38").Value = 85005
$ws.Range("C38").Value = 187
$ws.Range("D38").Value = 79608
$ws.Range("E38").Value = 4747
$ws.Range("G38").Value = 4
$ws.Range("H38").Value = 650

# Row 39
$ws.Range("A39").Value = "China"
$ws.Range("B39").Value = 85004
$ws.Range("C39").Value = 8
$ws.Range("D39").Value = 80046
$ws.Range("E39").Value = 324
$ws.Range("H39").Value = 4634

# Row 42
$ws.Range("B42").Value = 82945
$ws.Range("C42").Value = 674
$ws.Range("D42").Value = 74522
$ws.Range("E42").Value = 7901
$ws.Range("G42").Value = 1
$ws.Range("H42").Value = 522

# Row 43
$ws.Range("B43").Value = 71165
$ws.Range("C43").Value = 191
$ws.Range("D43").Value = 69650
$ws.Range("E43").Value = 853
$ws.Range("G43").Value = 5
$ws.Range("H43").Value = 662

# Row 45
$ws.Range("A45").Value = "Paises Bajos"
$ws.Range("B45").Value = 68624
$ws.Range("C45").Value = 510
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("G45").Value = 3
$ws.Range("H45").Value = 6218

# Row 46
$ws.Range("A46").Value = "Emiratos Arabes Unidos"
$ws.Range("B46").Value = 68511
$ws.Range("C46").Value = 491
$ws.Range("D46").Value = 59472
$ws.Range("E46").Value = 8661
$ws.Range("H46").Value = 378

# Row 68
$ws.Range("B68").Value = 33389
$ws.Range("C68").Value = 373
$ws.Range("D68").Value = 19368
$ws.Range("E68").Value = 13454
$ws.Range("G68").Value = 3
$ws.Range("H68").Value = 567

# Row 69
$ws.Range("B69").Value = 31099
$ws.Range("C69").Value = 125
$ws.Range("D69").Value = 29533
$ws.Range("E69").Value = 859

# Row 83
$ws.Range("B83").Value = 14592
$ws.Range("C83").Value = 38
$ws.Range("D83").Value = 13686
$ws.Range("E83").Value = 722
$ws.Range("G83").Value = 3
$ws.Range("H83").Value = 184

# Row 85
$ws.Range("B85").Value = 14004
$ws.Range("C85").Value = 90
$ws.Range("D85").Value = 10722
$ws.Range("E85").Value = 2697
$ws.Range("G85").Value = 7
$ws.Range("H85").Value = 585

# Row 142
$ws.Range("B142").Value = 1923
$ws.Range("C142").Value = 110
$ws.Range("D142").Value = 733
$ws.Range("E142").Value = 1142

# Row 185
$ws.Range("B185").Value = 272
$ws.Range("C185").Value = 2
$ws.Range("D185").Value = 216
$ws.Range("E185").Value = 56
